$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.028.33'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '1.644.77'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("D4").Value = "'" + '1.00'
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = "'" + '215.50'
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = "'" + '0.255'
$ws.Range("E8").Value = '  +0.28%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = "'" + '0.0639'
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("D10").Value = "'" + '19.60'
$ws.Range("E10").Value = '  -0.17%  '
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D12").Value = "'" + '4.26'
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("D13").Value = '1.637.97'
$ws.Range("E13").Value = '  +0.67%  '
$ws.Range("D14").Value = "'" + '0.544'
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("D15").Value = "'" + '63.46'
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").Value = '0.0₃0761'
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").Value = '26.073.20'
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").Value = "'" + '194.33'
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("E20").Value = '  -0.49%  '
$ws.Range("D21").Value = "'" + '9.91'
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("E22").Value = '  -1.03%  '
$ws.Range("D23").Value = "'" + '0.132'
$ws.Range("E23").Value = '  +4.71%  '
$ws.Range("D24").Value = "'" + '143.98'
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").Value = "'" + '6.89'
$ws.Range("E27").Value = '  +0.52%  '
$ws.Range("D28").Value = "'" + '15.50'
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("E29").Value = '  +0.33%  '
$ws.Range("E30").Value = '  -1.16%  '
$ws.Range("D31").Value = "'" + '3.26'
$ws.Range("E31").Value = '  +0.96%  '
$ws.Range("D32").Value = "'" + '3.28'
$ws.Range("E32").Value = '  -0.73%  '
$ws.Range("D33").Value = "'" + '1.54'
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("E34").Value = '  +1.07%  '
$ws.Range("D35").Value = "'" + '0.905'
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("D36").Value = '1.130.92'
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("D37").Value = "'" + '0.538'
$ws.Range("E37").Value = '  -1.41%  '
$ws.Range("D38").Value = "'" + '2.46'
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("D39").Value = "'" + '0.0157'
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("D41").Value = "'" + '98.95'
$ws.Range("E41").Value = '  -0.43%  '
$ws.Range("E42").Value = '  -0.30%  '
$ws.Range("E43").Value = '  +1.63%  '
$ws.Range("D44").Value = "'" + '56.50'
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").Value = "'" + '1.49'
$ws.Range("E45").Value = '  +2.39%  '
$ws.Range("E46").Value = '  -1.43%  '
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("D48").Value = "'" + '0.414'
$ws.Range("E48").Value = '  -0.27%  '
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("D50").Value = "'" + '0.0949'
$ws.Range("E50").Value = '  -1.36%  '
$ws.Range("D51").Value = "'" + '5.53'
$ws.Range("E51").Value = '  +0.04%  '
